# Apply commit "change Campaign to Topic" to advocacy_impact_dummy_data.xlsx
#
# Semantic changes:
#  - Sheet "campaign_mentions" renamed to "topic_mentions"
#  - Sheet "campaigns" renamed to "interventions"
#  - Column header "campaign" renamed to "topic" on every sheet that has it
#  - Column header "text" (on topic_mentions) renamed to "content"
#  - Data value "Car taxation" renamed to "Nature Restauration Law" wherever it
#    appears as a campaign/topic value
#  - View/selection cosmetic changes matching the authored diff

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename worksheets
# ---------------------------------------------------------------------------
$wsMentionsSupport = $wb.Worksheets.Item("mentions_support")
$wsTopicMentions    = $wb.Worksheets.Item("campaign_mentions")
$wsTopicMentions.Name = "topic_mentions"
$wsMepSentiment     = $wb.Worksheets.Item("mep_sentiment")
$wsInterventions    = $wb.Worksheets.Item("campaigns")
$wsInterventions.Name = "interventions"
$wsVotingResults    = $wb.Worksheets.Item("voting_results")

# ---------------------------------------------------------------------------
# 2. mentions_support: "campaign" header -> "topic", "Car taxation" -> "Nature Restauration Law"
# ---------------------------------------------------------------------------
$wsMentionsSupport.Range("C1").Value = "topic"
$wsMentionsSupport.Range("C2").Value = "Nature Restauration Law"

# ---------------------------------------------------------------------------
# 3. topic_mentions: "text" header -> "content", "campaign" header -> "topic",
#    "Car taxation" -> "Nature Restauration Law" for rows 2-11
# ---------------------------------------------------------------------------
$wsTopicMentions.Range("B1").Value = "content"
$wsTopicMentions.Range("D1").Value = "topic"
for ($r = 2; $r -le 11; $r++) {
    $wsTopicMentions.Cells.Item($r, 4).Value = "Nature Restauration Law"
}

# ---------------------------------------------------------------------------
# 4. mep_sentiment: "campaign" header -> "topic", "Car taxation" -> "Nature Restauration Law" for rows 2-13
# ---------------------------------------------------------------------------
$wsMepSentiment.Range("D1").Value = "topic"
for ($r = 2; $r -le 13; $r++) {
    $wsMepSentiment.Cells.Item($r, 4).Value = "Nature Restauration Law"
}

# ---------------------------------------------------------------------------
# 5. interventions (formerly campaigns): "campaign" header -> "topic",
#    "Car taxation" -> "Nature Restauration Law" for rows 2-11
# ---------------------------------------------------------------------------
$wsInterventions.Range("B1").Value = "topic"
for ($r = 2; $r -le 11; $r++) {
    $wsInterventions.Cells.Item($r, 2).Value = "Nature Restauration Law"
}

# ---------------------------------------------------------------------------
# 6. voting_results: "campaign" header -> "topic", "Car taxation" -> "Nature Restauration Law"
# ---------------------------------------------------------------------------
$wsVotingResults.Range("D1").Value = "topic"
$wsVotingResults.Range("D2").Value = "Nature Restauration Law"

# ---------------------------------------------------------------------------
# 7. View / selection cosmetic changes to match the authored diff
#    (select in file order so the final Activate leaves voting_results as the
#    tab-selected sheet, matching workbookView activeTab="4")
# ---------------------------------------------------------------------------
[void]$wsMentionsSupport.Range("C2").Select()
[void]$wsTopicMentions.Range("D9").Select()
[void]$wsMepSentiment.Range("K25").Select()
[void]$wsInterventions.Range("B2").Select()
[void]$wsVotingResults.Range("F24").Select()
